# API HW Error Fix
#
# Insert a blank paragraph and a new "AJAX: Asynchronous JavaScript and XML"
# paragraph immediately after the existing "JSON: JavaScript Object
# Notation ..." definition paragraph.

$d = $word.ActiveDocument

# Locate the paragraph that holds the JSON definition text.
$found = $d.Content
$null = $found.Find.Execute(
    "JSON: JavaScript Object Notation (Web format for sending data)",
    $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
$jsonPara = $found.Paragraphs(1)

# XML fragment for a blank "Times" paragraph (no run) -- matches the style
# already used for the blank lines elsewhere in this document.
$blankParaXml = '<w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:pPr><w:rPr><w:rFonts w:ascii="Times" w:hAnsi="Times"/></w:rPr></w:pPr></w:p>'

# Insert a blank paragraph right after the JSON paragraph.
$afterJson = $jsonPara.Next()
$insertionPoint1 = $d.Range($afterJson.Range.Start, $afterJson.Range.Start)
$null = $insertionPoint1.InsertXML($blankParaXml)

# Insert a second blank paragraph right after the one we just added; this
# one will receive the AJAX definition text below.
$afterBlank = $jsonPara.Next().Next()
$insertionPoint2 = $d.Range($afterBlank.Range.Start, $afterBlank.Range.Start)
$null = $insertionPoint2.InsertXML($blankParaXml)

# Fill the second new paragraph with the AJAX definition text, using the
# same "Times" font as the rest of the document.
$ajaxPara = $jsonPara.Next().Next()
$ajaxRange = $ajaxPara.Range
$null = $ajaxRange.MoveEnd(1, -1)
$ajaxRange.Text = "AJAX: Asynchronous JavaScript and XML"
$ajaxRange.Font.Name = "Times"
